# Update "want to go" counts (column F) for a handful of events that are
# duplicated across the "展览" (Exhibitions) sheet and the "全部类型"
# (All types) aggregate sheet.
#
#   F2: 4  -> 8
#   F6: 23 -> 24
#   F7: 5  -> 6

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 8
    $ws.Range("F6").Value = 24
    $ws.Range("F7").Value = 6
}
